$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$rows_PIR = @(
    @('2026-02-06', '10:20:21', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:20:25', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:20:28', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:20:31', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:20:35', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:20:40', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:20:45', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:20:50', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:20:55', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:21:00', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:21:05', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:21:10', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:21:15', '10:00', 'Bathroom', 'No Motion', 'Inactive'),
    @('2026-02-06', '10:21:20', '10:00', 'Bathroom', 'No Motion', 'Inactive')
)
$startRow_PIR = 486
$prefixCols_PIR = @(1)
$r = $startRow_PIR
foreach ($row in $rows_PIR) {
    for ($c = 1; $c -le 6; $c++) {
        $val = $row[$c - 1]
        if ($prefixCols_PIR -contains $c) {
            $val = "'" + $val
        }
        $ws.Cells.Item($r, $c).Value = $val
    }
    $r = $r + 1
}
Write-Host "PIR: wrote rows $startRow_PIR to $($r - 1)"

$ws = $wb.Worksheets.Item("Humidity")
$rows_Humidity = @(
    @('2026-02-06', '10:20:22', '10:00', 'Bathroom', '67.2%', 'Active'),
    @('2026-02-06', '10:20:26', '10:00', 'Bathroom', '68.2%', 'Active'),
    @('2026-02-06', '10:20:29', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:20:32', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:20:36', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:20:41', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:20:46', '10:00', 'Bathroom', '68.2%', 'Active'),
    @('2026-02-06', '10:20:51', '10:00', 'Bathroom', '68.2%', 'Active'),
    @('2026-02-06', '10:20:56', '10:00', 'Bathroom', '67.2%', 'Active'),
    @('2026-02-06', '10:21:01', '10:00', 'Bathroom', '68.2%', 'Active'),
    @('2026-02-06', '10:21:06', '10:00', 'Bathroom', '67.2%', 'Active'),
    @('2026-02-06', '10:21:11', '10:00', 'Bathroom', '68.3%', 'Active'),
    @('2026-02-06', '10:21:16', '10:00', 'Bathroom', '67.2%', 'Active')
)
$startRow_Humidity = 337
$prefixCols_Humidity = @(1, 5)
$r = $startRow_Humidity
foreach ($row in $rows_Humidity) {
    for ($c = 1; $c -le 6; $c++) {
        $val = $row[$c - 1]
        if ($prefixCols_Humidity -contains $c) {
            $val = "'" + $val
        }
        $ws.Cells.Item($r, $c).Value = $val
    }
    $r = $r + 1
}
Write-Host "Humidity: wrote rows $startRow_Humidity to $($r - 1)"

$ws = $wb.Worksheets.Item("Temperature")
$rows_Temperature = @(
    @('2026-02-06', '10:20:24', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:20:27', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:20:30', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:20:33', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:20:37', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:20:42', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:20:47', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:20:52', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:20:57', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:21:02', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:21:07', '10:00', 'Bathroom', '28.2C', 'Active'),
    @('2026-02-06', '10:21:12', '10:00', 'Bathroom', '28.3C', 'Active'),
    @('2026-02-06', '10:21:17', '10:00', 'Bathroom', '28.2C', 'Active')
)
$startRow_Temperature = 337
$prefixCols_Temperature = @(1)
$r = $startRow_Temperature
foreach ($row in $rows_Temperature) {
    for ($c = 1; $c -le 6; $c++) {
        $val = $row[$c - 1]
        if ($prefixCols_Temperature -contains $c) {
            $val = "'" + $val
        }
        $ws.Cells.Item($r, $c).Value = $val
    }
    $r = $r + 1
}
Write-Host "Temperature: wrote rows $startRow_Temperature to $($r - 1)"
